$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'258.70"
$ws.Range("E2").Value = "'5.19%"

$ws.Range("D3").Value = "'27.28"
$ws.Range("E3").Value = "'-3.63%"

$ws.Range("D4").Value = "'5.220"
$ws.Range("E4").Value = "'-1.27%"

$ws.Range("D5").Value = "'0.05932"
$ws.Range("E5").Value = "'3.86%"

$ws.Range("D6").Value = "'6.713"
$ws.Range("E6").Value = "'1.11%"

$ws.Range("D7").Value = "'0.8657"
$ws.Range("E7").Value = "'0.36%"

$ws.Range("D8").Value = "'0.9991"
$ws.Range("E8").Value = "'12.06%"

$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1418"
$ws.Range("E9").Value = "'2.17%"

$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.07165"
$ws.Range("E10").Value = "'1.22%"

$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.03143"
$ws.Range("E11").Value = "'-0.32%"

$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.09226"
$ws.Range("E12").Value = "'-0.04%"

$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001546"
$ws.Range("E13").Value = "'1.38%"

$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").Value = "'0.0006047"
$ws.Range("E14").Value = "'-94.06%"

$ws.Range("D15").Value = "'0.005892"
$ws.Range("E15").Value = "'-3.25%"

$ws.Range("E16").Value = "'0.09%"

$ws.Range("D17").Value = "'3.269"
$ws.Range("E17").Value = "'1.70%"

$ws.Range("D18").Value = "'2.227"
$ws.Range("E18").Value = "'2.49%"

$ws.Range("E19").Value = "'-0.70%"

$ws.Range("D20").Value = "'0.03553"
$ws.Range("E20").Value = "'6.90%"

$ws.Range("E21").Value = "'-0.21%"

$ws.Range("D22").Value = "'3.536"
$ws.Range("E22").Value = "'1.61%"

$ws.Range("D23").Value = "'0.04184"
$ws.Range("E23").Value = "'2.27%"

$ws.Range("E24").Value = "'1.48%"

$ws.Range("D25").Value = "'0.001216"
$ws.Range("E25").Value = "'-0.36%"

$ws.Range("D26").Value = "'0.004525"

$ws.Range("D27").Value = "'0.0001198"
$ws.Range("E27").Value = "'-0.14%"

$ws.Range("D28").Value = "'0.0001937"
$ws.Range("E28").Value = "'34.06%"

$ws.Range("D40").Value = "'0.03837"
$ws.Range("E40").Value = "'1.14%"

$ws.Range("D41").Value = "'0.006570"
$ws.Range("E41").Value = "'16.30%"

$ws.Range("D42").Value = "'0.1103"

$ws.Range("D43").Value = "'0.002276"
$ws.Range("E43").Value = "'3.46%"

$ws.Range("D44").Value = "'0.01075"
$ws.Range("E44").Value = "'13.24%"

$ws.Range("D45").Value = "'0.00005433"
$ws.Range("E45").Value = "'2.79%"

$ws.Range("E46").Value = "'-0.13%"

$ws.Range("D47").Value = "'0.1091"
$ws.Range("E47").Value = "'22.41%"

$ws.Range("D48").Value = "'0.002230"
$ws.Range("E48").Value = "'-1.42%"

$ws.Range("D49").Value = "'0.00002097"
$ws.Range("E49").Value = "'-0.13%"

$ws.Range("D50").Value = "'0.0001997"
$ws.Range("E50").Value = "'-0.13%"
